# Daily attendance processing - 2025-10-26 18:26:21
# Reorders the comma-separated "Recorded By" names/emails in column G so
# that "System" (and similar already-sorted duplicate entries) is moved to
# the end of the list instead of the front.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact before -> after string replacements observed for the "Recorded By"
# (column G) values across the sheet.
$map = @{
    "backup@backdoor.com, System, system" = "backup@backdoor.com, system, System"
    "System, dnasr281@gmail.com"           = "dnasr281@gmail.com, System"
    "System, admin@admin.com"              = "admin@admin.com, System"
    "admin@admin.com, dnasr281@gmail.com"  = "dnasr281@gmail.com, admin@admin.com"
}

# Determine the last used row in the sheet so we cover every data row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = $ws.UsedRange.Rows.Count }

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value = $map[$val]
    }
}
